$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $range = $d.Content
    $found = $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $found) {
        throw "Could not find text: $findText"
    }
}

# Title (Heading1) and the bold "title" run near the bottom of the document
# Both occurrences share identical text; wdReplaceAll (2) replaces them both in one call.
Replace-Text "Play Mistress of Egypt for Free - Game Review & Ratings" "Play Mistress of Egypt Free: Exciting Gameplay and Huge Jackpots"

# "What we like" bullet list
Replace-Text "Three progressive jackpots" "Synchronized reels"
Replace-Text "Mesmerizing Middle Eastern soundtrack" "Random wilds"
Replace-Text "Synchronized reels and random wilds" "Free spin round with high potential wins"
Replace-Text "Wins of up to 200x bet on each spin" "Mesmerizing graphics and soundtrack"

# "What we don't like" bullet list
Replace-Text "Medium variance gameplay" "Medium variance"
Replace-Text "Base RTP of 95.60%" "Base RTP is slightly low"

# Meta description (italic run at the end of the document)
Replace-Text "Discover the ancient Egypt-themed online slot Mistress of Egypt. Play for free with synced reels, random wilds, & three progressive jackpots. Read our full review." "Play Mistress of Egypt free and enjoy synchronized reels, random wilds, and big wins!"
